# Add a new interview-experience record (row 11) to the log, mirroring the
# formatting of the previous rows (date style, wrapped-question style, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the prior row (row 10) down into the new row so the
# new cells pick up the same cell styles (date format, borders, wrap text).
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)  # xlPasteFormats

# Date: 20-May-2025 (serial 45797)
$ws.Range("A11").Value = 45797

# Company / round
$ws.Range("B11").Value = "cresensolutions - technical round"

# Interview questions asked
$ws.Range("C11").Value = "string isblank vs isempty, design pattern - factory pattern, singletone pattern , es and redis use, hashmap and hashcollision, final , finally , finalize, handle multiple exception in single catch, metaannotaion in springboot, implement auth using properties file before requesting to controller,  shallow copy and deep copy"

# Result
$ws.Range("D11").Value = "cleared"

# Match the taller row height Excel produced for the wrapped question text.
$ws.Rows.Item(11).RowHeight = 90

# Reflect the user's final selection/cursor position after the edit.
$ws.Range("C11").Select()
